$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The activity stats (runs/balls/fours/sixes columns C:F) for rows 2-13
# have been refreshed/reordered to reflect the player's updated match log.
# Ensure the cells stay text-typed (matching the source workbook's
# "numberStoredAsText" cells) before writing the new values.
$ws.Range("C2:F13").NumberFormat = "@"

$data = @{
    2  = @("34", "30", "4", "1")
    3  = @("0",  "3",  "0", "0")
    4  = @("10", "10", "1", "0")
    5  = @("3",  "4",  "0", "0")
    6  = @("13", "10", "0", "1")
    7  = @("12", "17", "0", "1")
    8  = @("12", "7",  "2", "0")
    9  = @("3",  "3",  "0", "0")
    10 = @("1",  "1",  "0", "0")
    11 = @("1",  "3",  "0", "0")
    12 = @("20", "4",  "2", "2")
    13 = @("0",  "0",  "0", "0")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("D$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
}
